$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values for rows 2-8 (names unchanged)
$ws.Range("B2").Value = 37
$ws.Range("B3").Value = 32
$ws.Range("B4").Value = 29
$ws.Range("B5").Value = 26
$ws.Range("B6").Value = 26
$ws.Range("B7").Value = 24
$ws.Range("B8").Value = 24

# Rows 9-11: names get reordered, and B9/B10 values change (B11 stays the same)
$ws.Range("A9").Value = "URBINA ANTICONA ALEX BRUNO"
$ws.Range("B9").Value = 22

$ws.Range("A10").Value = "LLANOS HUACCHA BRITSY"
$ws.Range("B10").Value = 20

$ws.Range("A11").Value = "OLIVA MUÑOZ LOURDES ARACELY"
